# Weekly update: insert two new daily price records for Coliflor /
# Vega Modelo de Temuco, right before the existing row 449 block.
# This pushes the existing rows 449:546 down to 451:548 (matching the
# new sheet dimension A1:R548) while keeping all of their data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("449:450").Insert()

# New row 449 - "Primera" quality record
$ws.Cells.Item(449, 1).Value  = 10
$ws.Cells.Item(449, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(449, 3).Value  = "La Araucanía"
$ws.Cells.Item(449, 4).Value  = 44964
$ws.Cells.Item(449, 5).Value  = 9
$ws.Cells.Item(449, 6).Value  = 100112008
$ws.Cells.Item(449, 7).Value  = "Coliflor"
$ws.Cells.Item(449, 8).Value  = "Sin especificar"
$ws.Cells.Item(449, 9).Value  = "Primera"
$ws.Cells.Item(449, 10).Value = 300
$ws.Cells.Item(449, 11).Value = 1300
$ws.Cells.Item(449, 12).Value = 1300
$ws.Cells.Item(449, 13).Value = 1300
$ws.Cells.Item(449, 14).Value = "`$/unidad"
$ws.Cells.Item(449, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(449, 16).Value = 1300
$ws.Cells.Item(449, 17).Value = 1
$ws.Cells.Item(449, 18).Value = "Hortaliza"

# New row 450 - "Segunda" quality record
$ws.Cells.Item(450, 1).Value  = 10
$ws.Cells.Item(450, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(450, 3).Value  = "La Araucanía"
$ws.Cells.Item(450, 4).Value  = 44964
$ws.Cells.Item(450, 5).Value  = 9
$ws.Cells.Item(450, 6).Value  = 100112008
$ws.Cells.Item(450, 7).Value  = "Coliflor"
$ws.Cells.Item(450, 8).Value  = "Sin especificar"
$ws.Cells.Item(450, 9).Value  = "Segunda"
$ws.Cells.Item(450, 10).Value = 300
$ws.Cells.Item(450, 11).Value = 1200
$ws.Cells.Item(450, 12).Value = 1200
$ws.Cells.Item(450, 13).Value = 1200
$ws.Cells.Item(450, 14).Value = "`$/unidad"
$ws.Cells.Item(450, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(450, 16).Value = 1200
$ws.Cells.Item(450, 17).Value = 1
$ws.Cells.Item(450, 18).Value = "Hortaliza"
